$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns G and H: TexturePacker / UIAtlas support for thumbnail image rendering
# Row 1 - field names
$ws.Range("G1").Value = "simple_atlas"
$ws.Range("H1").Value = "simple_tex"

# Row 2 - field types
$ws.Range("G2").Value = "string"
$ws.Range("H2").Value = "string"

# Row 3 - Chinese display labels
$ws.Range("G3").Value = "缩略图集"
$ws.Range("H3").Value = "缩略图"

# Row 4 - sample data row
$ws.Range("G4").Value = "CardSimple"

# Move/extend selection the same way the authoring app left it
$ws.Range("H5").Select()
